$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("time_estimation_sheet")

$ws.Range("E10").Value = "2hrs"
$ws.Range("E11").Value = "2 hrs"
$ws.Range("E12").Value = "1 hr"
$ws.Range("E13").Value = "15 min"
$ws.Range("E14").Value = "1hr"
$ws.Range("E15").Value = "1hr"
$ws.Range("E16").Value = "2hr"
$ws.Range("E17").Value = "1hr"
$ws.Range("E19").Value = "1hr"
$ws.Range("E20").Value = "2hr"
